# Weekly update: insert two new price observations for "Zapallo italiano"
# right before the existing row 379 (old data), shifting subsequent rows
# down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 379 (rows 379 & 380 become new rows;
# everything that used to start at row 379 shifts down to row 381).
$ws.Rows("379:380").Insert()

# Columns A, B, C, E, F, G, R are constant across this whole data block,
# so copy them from the (unchanged) row directly above the insertion.
$constCols = @(1, 2, 3, 5, 6, 7, 18)
foreach ($col in $constCols) {
    $srcVal = $ws.Cells.Item(378, $col).Value2
    $ws.Cells.Item(379, $col).Value = $srcVal
    $ws.Cells.Item(380, $col).Value = $srcVal
}

# --- New row 379 ---
$ws.Cells.Item(379, 4).Value  = 44714                                  # D: Fecha
$ws.Cells.Item(379, 8).Value  = 'Bola 8'                                # H: Variedad
$ws.Cells.Item(379, 9).Value  = 'Primera'                               # I: Calidad
$ws.Cells.Item(379, 10).Value = 400                                     # J: Volumen
$ws.Cells.Item(379, 11).Value = 8000                                    # K: Precio minimo
$ws.Cells.Item(379, 12).Value = 10000                                   # L: Precio maximo
$ws.Cells.Item(379, 13).Value = 9150                                    # M: Precio promedio ponderado
$ws.Cells.Item(379, 14).Value = '$/caja 50 unidades'                    # N: Unidad de comercializacion
$ws.Cells.Item(379, 15).Value = 'Región de Arica y Parinacota'          # O: Origen
$ws.Cells.Item(379, 16).Value = 183                                     # P: Precio $/Kg
$ws.Cells.Item(379, 17).Value = 50                                      # Q: Kg o Unidades

# --- New row 380 ---
$ws.Cells.Item(380, 4).Value  = 44714                                   # D: Fecha
$ws.Cells.Item(380, 8).Value  = 'Sin especificar'                       # H: Variedad
$ws.Cells.Item(380, 9).Value  = 'Primera'                               # I: Calidad
$ws.Cells.Item(380, 10).Value = 740                                     # J: Volumen
$ws.Cells.Item(380, 11).Value = 10000                                   # K: Precio minimo
$ws.Cells.Item(380, 12).Value = 12000                                   # L: Precio maximo
$ws.Cells.Item(380, 13).Value = 11135                                   # M: Precio promedio ponderado
$ws.Cells.Item(380, 14).Value = '$/caja 50 unidades'                    # N: Unidad de comercializacion
$ws.Cells.Item(380, 15).Value = 'Región de Arica y Parinacota'          # O: Origen
$ws.Cells.Item(380, 16).Value = 223                                     # P: Precio $/Kg
$ws.Cells.Item(380, 17).Value = 50                                      # Q: Kg o Unidades
